$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "AutoCustoPKxw_0105599"
$ws.Range("D2").Value = "zBpiTTUFDf"
$ws.Range("C3").Value = "AutoCustvjMeg_0105972"
$ws.Range("D3").Value = "YIRBCPgopR"
